$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2 through 9
# from 45204 (2023-10-05) to 45207 (2023-10-08), keeping existing styling.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45207
}
